# Populate the Clothing column (G) for rows 2-86 based on the trained clothing model output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
    2 = 'Trunks,Jodhpurs'
    3 = 'Jumpsuit,Kaftan'
    4 = 'Caftan,Trunks'
    5 = 'Jumpsuit,Blouse'
    6 = 'Jumpsuit,Kaftan'
    7 = 'Jumpsuit,Blouse'
    8 = 'Jumpsuit,Blouse'
    9 = 'Caftan,Jodhpurs'
    10 = 'Jumpsuit,Halter'
    11 = 'Trunks,Blazer'
    12 = 'Jumpsuit,Halter'
    13 = 'Trunks,Blazer'
    14 = 'Jumpsuit,Blouse'
    15 = 'Jumpsuit,Kaftan'
    16 = 'Caftan,Parka'
    17 = 'Jumpsuit,Trunks'
    18 = 'Halter,Jumpsuit'
    19 = 'Jumpsuit,Blazer'
    20 = 'Blazer,Top'
    21 = 'Jumpsuit,Turtleneck'
    22 = 'Halter,Caftan'
    23 = 'Trunks,Turtleneck'
    24 = 'Jumpsuit,Kaftan'
    25 = 'Jumpsuit,Trunks'
    26 = 'Jumpsuit,Kaftan'
    27 = 'Jumpsuit,Blouse'
    28 = 'Caftan,Parka'
    29 = 'Caftan,Blazer'
    30 = 'Trunks,Jodhpurs'
    31 = 'Jumpsuit,Kaftan'
    32 = 'Jumpsuit,Blouse'
    33 = 'Halter,Jumpsuit'
    34 = 'Jumpsuit,Dress'
    35 = 'Kaftan,Halter'
    36 = 'Caftan,Dress'
    37 = 'Halter,Blazer'
    38 = 'Trunks,Jumpsuit'
    39 = 'Parka,Dress'
    40 = 'Trunks,Jodhpurs'
    41 = 'Halter,Blazer'
    42 = 'Halter,Tee'
    43 = 'Jumpsuit,Dress'
    44 = 'Jumpsuit,Blouse'
    45 = 'Halter,Blazer'
    46 = 'Parka,Dress'
    47 = 'Parka,Jumpsuit'
    48 = 'Jumpsuit,Dress'
    49 = 'Caftan,Blouse'
    50 = 'Parka,Jumpsuit'
    51 = 'Blouse,Jumpsuit'
    52 = 'Jumpsuit,Trunks'
    53 = 'Jumpsuit,Kaftan'
    54 = 'Jumpsuit,Trunks'
    55 = 'Dress,Jumpsuit'
    56 = 'Jumpsuit,Kaftan'
    57 = 'Jumpsuit,Halter'
    58 = 'Caftan,Parka'
    59 = 'Jumpsuit,Jodhpurs'
    60 = 'Caftan,Trunks'
    61 = 'Dress,Jumpsuit'
    62 = 'Halter,Trunks'
    63 = 'Dress,Jumpsuit'
    64 = 'Jumpsuit,Kaftan'
    65 = 'Blouse,Parka'
    66 = 'Jodhpurs,Capris'
    67 = 'Jumpsuit,Dress'
    68 = 'Halter,Tee'
    69 = 'Jumpsuit,Blouse'
    70 = 'Jumpsuit,Halter'
    71 = 'Jumpsuit,Kaftan'
    72 = 'Jumpsuit,Kaftan'
    73 = 'Caftan,Trunks'
    74 = 'Jodhpurs,Trunks'
    75 = 'Jumpsuit,Blouse'
    76 = 'Blazer,Trunks'
    77 = 'Parka,Gauchos'
    78 = 'Caftan,Jumpsuit'
    79 = 'Trunks,Caftan'
    80 = 'Jumpsuit,Blouse'
    81 = 'Jumpsuit,Dress'
    82 = 'Blazer,Trunks'
    83 = 'Trunks,Caftan'
    84 = 'Blazer,Trunks'
    85 = 'Caftan,Halter'
    86 = 'Caftan,Trunks'
}

foreach ($row in $clothing.Keys) {
    $ws.Cells.Item($row, 7).Value = $clothing[$row]
}
